$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 730
$ws.Range("AO2").Value = 580
$ws.Range("AM3").Value = 95
$ws.Range("R4").Value = 1.81
$ws.Range("S4").Value = 2.2
$ws.Range("F5").Value = 2.2
$ws.Range("G5").Value = 2.22
$ws.Range("K5").Value = 3.85
$ws.Range("L5").Value = 1.32
$ws.Range("W5").Value = 1.82
$ws.Range("Q6").Value = 2.18
$ws.Range("P7").Value = 1.97
$ws.Range("Q7").Value = 2
$ws.Range("T7").Value = 1.78
$ws.Range("U7").Value = 2.24
$ws.Range("F8").Value = 1.91
$ws.Range("G8").Value = 1.92
$ws.Range("H8").Value = 4.7
$ws.Range("I8").Value = 4.8
$ws.Range("V8").Value = 1.26
$ws.Range("W8").Value = 2.08
$ws.Range("X8").Value = 12.5
$ws.Range("Z8").Value = 34
$ws.Range("AA8").Value = 110
$ws.Range("AD8").Value = 18.5
$ws.Range("H9").Value = 2.14
$ws.Range("I9").Value = 2.16
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 3.8
$ws.Range("V9").Value = 1.86
$ws.Range("X9").Value = 16.5
$ws.Range("F10").Value = 2.56
$ws.Range("G10").Value = 2.58
$ws.Range("L10").Value = 1.46
$ws.Range("V10").Value = 1.42
$ws.Range("W10").Value = 1.63
$ws.Range("Y10").Value = 10.5
$ws.Range("F11").Value = 2.4
$ws.Range("G11").Value = 2.42
$ws.Range("S11").Value = 4.8
$ws.Range("U11").Value = 1.93
$ws.Range("W11").Value = 1.7
$ws.Range("L12").Value = 1.28
$ws.Range("U12").Value = 1.96
$ws.Range("AK12").Value = 150
$ws.Range("AN12").Value = 160
$ws.Range("L13").Value = 1.37
$ws.Range("T13").Value = 1.85
$ws.Range("AC13").Value = 8.8
$ws.Range("F14").Value = 5.4
$ws.Range("G14").Value = 5.5
$ws.Range("L14").Value = 1.34
$ws.Range("N14").Value = 4.4
$ws.Range("T14").Value = 1.79
$ws.Range("V14").Value = 2.32
$ws.Range("AC14").Value = 8.8
$ws.Range("AM14").Value = 90
